# Update scripts with new TPM: recomputed Receptor-expressing-cells and all
# downstream derived columns for the Fgf15-Fgfr4 LR-pair sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs) -- Receptor-expressing cells count changed 1 -> 2
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.08378199999999998
$ws.Range("N2").Value = 0.251346
$ws.Range("O2").Value = 0.007571394704126512
$ws.Range("P2").Value = 0.007571394704126512
$ws.Range("Q2").Value = 0.0112225989
$ws.Range("R2").Value = 0.1010033901
$ws.Range("S2").Value = 0.007571394704126512
$ws.Range("T2").Value = 0.007571394704126512

# Row 3 (Target cluster: FAPs) -- minor re-normalization from the row 2 update
$ws.Range("M3").Value = 0.07352966666666667
$ws.Range("O3").Value = 0.006644889460697858
$ws.Range("P3").Value = 0.006644889460697857
$ws.Range("Q3").Value = 0.009849298850000002
$ws.Range("R3").Value = 0.08864368965000001
$ws.Range("S3").Value = 0.006644889460697858
$ws.Range("T3").Value = 0.006644889460697857

# Row 4 (Target cluster: MuSCs) -- minor re-normalization from the row 2 update
$ws.Range("O4").Value = 0.9857837158351757
$ws.Range("P4").Value = 0.9857837158351755
$ws.Range("S4").Value = 0.9857837158351757
$ws.Range("T4").Value = 0.9857837158351755
